$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# Rebuild the "Email/password" report into a 4-column "Test Case" report.
# Old col A (Email/links) and col B (password/numbers) move right to C/D;
# new col A holds the run flag, new col B holds the test-case description.
#
# Values are written in this specific order so the regenerated shared-
# string table lands in the same order as the target workbook:
#   0 Email, 1 password, 2 wiasm.mtour@gmail.com, 3 Test Case, 4 Run,
#   5 Check response..., 6 asdas, 7 asd, 8 retyertetertert, 9 rterterter
# -----------------------------------------------------------------------

$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "password"
$ws.Range("B1").Value = "Test Case"
$ws.Range("A1").Value = "Run"

$ws.Range("C2").Value = "wiasm.mtour@gmail.com"
$ws.Range("D2").Value = 123456789
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Check response on entering valid  Credentials(Email and password)"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "asdas"
$ws.Range("C3").Value = "asd"
$ws.Range("D3").Value = "asd"

$ws.Range("C4").Value = "wiasm.mtour@gmail.com"
$ws.Range("D4").Value = 123456789
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Check response on entering valid  Credentials(Email and password)"

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "retyertetertert"
$ws.Range("C5").Value = "rterterter"
$ws.Range("D5").Value = 25121

# -----------------------------------------------------------------------
# Hyperlinks: drop the four stale mailto links and recreate just the two
# that are still needed, now anchored on column C. Do this before the
# style fix-up below, since Hyperlinks.Add() stamps its own font style
# on the target cell that we don't want to keep.
# -----------------------------------------------------------------------

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:wiasm.mtour@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:wiasm.mtour@gmail.com") | Out-Null

# -----------------------------------------------------------------------
# Styles: reuse the two existing cell formats (plain left-aligned text,
# and the underlined "Hyperlink" text) instead of creating new ones.
# -----------------------------------------------------------------------

# Plain left-aligned style (already used by A1/B1) -> new/changed cells
# that should NOT look like a hyperlink.
$ws.Range("A1").Copy() | Out-Null
foreach ($addr in @("C1", "D1", "D2", "A3", "C3", "D3", "D4", "A5", "C5", "D5")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Underlined "Hyperlink" style (already used by A2/A3/A4/A5) -> cells that
# should look like a hyperlink / belong to the highlighted rows.
$ws.Range("A2").Copy() | Out-Null
foreach ($addr in @("B2", "C2", "B4", "C4")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# -----------------------------------------------------------------------
# Column widths for the new/changed columns B, C, D (A keeps its width).
# -----------------------------------------------------------------------

$ws.Columns("B").ColumnWidth = 62.43
$ws.Columns("C").ColumnWidth = 30.26
$ws.Columns("D").ColumnWidth = 37.6

# -----------------------------------------------------------------------
# Page setup / selection to match the refreshed sheet.
# -----------------------------------------------------------------------

$ws.PageSetup.Orientation = 1
$ws.Range("A6:XFD36").Select()
